$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("B2").Value = "STEELLEVELINMOLD-SETMEAN"
$ws.Range("C2").Value = 0.02484936562097117

# Update row 3
$ws.Range("B3").Value = "NARROWFACESECONDARYCOOLINGWATERPRESSUREACTUALMEAN"
$ws.Range("C3").Value = 0.02430368112155793

# Update row 4
$ws.Range("A4").Value = 3
$ws.Range("C4").Value = 0.01622696814615959

# Update row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "TUNDISHWEIGHTMEAN"
$ws.Range("C5").Value = 0.01533290323001619

# Update row 6 (value only)
$ws.Range("C6").Value = 0.01415176609373432

# Update row 7
$ws.Range("B7").Value = "NARROWFACE1WATERFLOWMEAN"
$ws.Range("C7").Value = 0.01081853183931805

# Update row 8
$ws.Range("A8").Value = 13
$ws.Range("B8").Value = "SEG9+10+11+12BOTTOMSECONDARYCOOLINGWATERFLOWACTUALMEAN"
$ws.Range("C8").Value = 0.009801561493585765

# Update row 9
$ws.Range("A9").Value = 14
$ws.Range("B9").Value = "LIQUIDUSTEMPMEAN"
$ws.Range("C9").Value = 0.009688317844439423
